$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

# Fixed, never-modified cells used purely as a style source so that newly
# populated cells pick up the same formatting as the rest of the header /
# data rows (B1 = bold header style, B2 = plain data-row style).
$wsStyleSrc = $wb.Worksheets.Item("土地")
$headerStyleCell = $wsStyleSrc.Range("B1")
$dataStyleCell   = $wsStyleSrc.Range("B2")

# Set a text value on a cell in a way that never gets auto-converted to a
# number/date by Excel's input parser, then stamp it with the same style as
# a reference cell (so new cells match the look of the existing header / data
# rows instead of picking up stray formatting).
function Set-TextCell($ws, $addr, $text, $isHeader) {
    $ws.Range($addr).Value = "'" + $text
    if ($isHeader) {
        $headerStyleCell.Copy() | Out-Null
    } else {
        $dataStyleCell.Copy() | Out-Null
    }
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
}

# Set a numeric value on a cell and stamp it with the same style as a
# reference cell.
function Set-NumberCell($ws, $addr, $number, $isHeader) {
    if ($isHeader) {
        $headerStyleCell.Copy() | Out-Null
    } else {
        $dataStyleCell.Copy() | Out-Null
    }
    $ws.Range($addr).PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Range($addr).Value = $number
}

# ---------------------------------------------------------------
# Sheet "保險" (insurance) -- sheet3
# ---------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# Header row: relabel existing header cells, add new header cells E1:K1
Set-TextCell $wsIns "B1" "company"           $true
Set-TextCell $wsIns "C1" "name"              $true
Set-TextCell $wsIns "D1" "owner"             $true
Set-TextCell $wsIns "E1" "property_category" $true
Set-TextCell $wsIns "F1" "category"          $true
Set-TextCell $wsIns "G1" "date"              $true
Set-TextCell $wsIns "H1" "legislator_name"   $true
Set-TextCell $wsIns "I1" "legislator_id"     $true
Set-TextCell $wsIns "J1" "source_file"       $true
Set-TextCell $wsIns "K1" "index"             $true

# Data row 2 (index 101): B2/C2/D2 already correct, add E2:K2
Set-TextCell   $wsIns "E2" "insurance"   $false
Set-TextCell   $wsIns "F2" "normal"      $false
Set-TextCell   $wsIns "G2" "2012-04-20"  $false
Set-TextCell   $wsIns "H2" "陳超明"       $false
Set-NumberCell $wsIns "I2" 836           $false
Set-TextCell   $wsIns "J2" "tmpb48f1"    $false
Set-NumberCell $wsIns "K2" 101           $false

# Data row 3 (index 102): B3/C3/D3 already correct, add E3:K3
Set-TextCell   $wsIns "E3" "insurance"   $false
Set-TextCell   $wsIns "F3" "normal"      $false
Set-TextCell   $wsIns "G3" "2012-04-20"  $false
Set-TextCell   $wsIns "H3" "陳超明"       $false
Set-NumberCell $wsIns "I3" 836           $false
Set-TextCell   $wsIns "J3" "tmpb48f1"    $false
Set-NumberCell $wsIns "K3" 102           $false

# ---------------------------------------------------------------
# Sheet "債務" (debt) -- sheet4
# ---------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item("債務")

# Header row: relabel existing header cells, add new header cells H1:N1
Set-TextCell $wsDebt "B1" "species"           $true
Set-TextCell $wsDebt "C1" "debtor"            $true
Set-TextCell $wsDebt "D1" "owner"             $true
Set-TextCell $wsDebt "E1" "total"             $true
Set-TextCell $wsDebt "F1" "register_date"     $true
Set-TextCell $wsDebt "G1" "register_reason"   $true
Set-TextCell $wsDebt "H1" "property_category" $true
Set-TextCell $wsDebt "I1" "category"          $true
Set-TextCell $wsDebt "J1" "date"              $true
Set-TextCell $wsDebt "K1" "legislator_name"   $true
Set-TextCell $wsDebt "L1" "legislator_id"     $true
Set-TextCell $wsDebt "M1" "source_file"       $true
Set-TextCell $wsDebt "N1" "index"             $true

# Data row 2 (index 113): B2:G2 already correct, add H2:N2
Set-TextCell   $wsDebt "H2" "debt"        $false
Set-TextCell   $wsDebt "I2" "normal"      $false
Set-TextCell   $wsDebt "J2" "2012-04-20"  $false
Set-TextCell   $wsDebt "K2" "陳超明"       $false
Set-NumberCell $wsDebt "L2" 836           $false
Set-TextCell   $wsDebt "M2" "tmpb48f1"    $false
Set-NumberCell $wsDebt "N2" 113           $false

# ---------------------------------------------------------------
# Sheet "事業投資" (investment) -- sheet5
# ---------------------------------------------------------------
$wsInv = $wb.Worksheets.Item("事業投資")

# Header row: relabel existing header cells, add new header cells H1:N1
Set-TextCell $wsInv "B1" "owner"             $true
Set-TextCell $wsInv "C1" "company"           $true
Set-TextCell $wsInv "D1" "address"           $true
Set-TextCell $wsInv "E1" "total"             $true
Set-TextCell $wsInv "F1" "register_date"     $true
Set-TextCell $wsInv "G1" "register_reason"   $true
Set-TextCell $wsInv "H1" "property_category" $true
Set-TextCell $wsInv "I1" "category"          $true
Set-TextCell $wsInv "J1" "date"              $true
Set-TextCell $wsInv "K1" "legislator_name"   $true
Set-TextCell $wsInv "L1" "legislator_id"     $true
Set-TextCell $wsInv "M1" "source_file"       $true
Set-TextCell $wsInv "N1" "index"             $true

# Data row 2 (index 119): B2:G2 already correct, add H2:N2
Set-TextCell   $wsInv "H2" "investment"  $false
Set-TextCell   $wsInv "I2" "normal"      $false
Set-TextCell   $wsInv "J2" "2012-04-20"  $false
Set-TextCell   $wsInv "K2" "陳超明"       $false
Set-NumberCell $wsInv "L2" 836           $false
Set-TextCell   $wsInv "M2" "tmpb48f1"    $false
Set-NumberCell $wsInv "N2" 119           $false
